# Apply the "ProjectsSubTypes" parameter-sheet update:
#  - employee quarterly rates
#  - names of functional groups
#  - hourly production calendar
#  - other new reference values
#
# Strategy: first reshape the sheet (insert/delete whole rows so the
# row count / blank-separator layout matches the target), then write the
# final text for every data cell directly so content is exactly right
# regardless of how Excel defaulted newly inserted cells.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# ---------------------------------------------------------------------
# Step 1: structural row inserts (bottom-up, so earlier row numbers used
# below stay valid while we work).
# ---------------------------------------------------------------------
$ws.Rows("114").Insert()
$ws.Rows("95").Insert()
$ws.Rows("59").Insert()
$ws.Rows("56").Insert()
$ws.Rows("43:44").Insert()
$ws.Rows(45).Delete()
$ws.Rows("38").Insert()
$ws.Rows("35:36").Insert()
$ws.Rows("26:27").Insert()
$ws.Rows("14:16").Insert()

# ---------------------------------------------------------------------
# Step 2: write final values for every row (1-based), column A then B.
# An empty string clears the cell back to a blank (but still styled) cell.
# ---------------------------------------------------------------------
$rows = @(
  @("ProjectName", "ProjectSubTypePart"),
  @("А0172. Адм.время НУП. Отпуска", "От"),
  @("А0212-Адм. время НИСУ. Отпуск", "От"),
  @("А0214-Адм. время НКиПС. Отпуск", "От"),
  @("А0218-Адм. время НСКК. Отпуск", "От"),
  @("А0222-Адм. время НПИС. Отпуск", "От"),
  @("А0223-Адм. время НПУ. Отпуск", "От"),
  @("А0224-Адм. время НР. Отпуск", "От"),
  @("А0295 -Адм. Время. НРБ. Отпуск", "От"),
  @("А0428 - Отгулы за`u{202f}ранее`u{202f}отработанное время (переработки)", "От"),
  @("А0524 - Административное время ООД. Отпуск.", "От"),
  @("А0425 - Отгулы за`u{202f}ранее`u{202f}отработанное время (переработки)", "От"),
  @("А0426 - Отгулы за`u{202f}ранее`u{202f}отработанное время (переработки)", "От"),
  @("А0602 - Адм. Время. RCBI. Отпуск", "От"),
  @("А0212-Адм. время ДИСУ. Отпуск", "От"),
  @("", ""),
  @("", ""),
  @("А0523 - Административное время ООД. Управление персоналом.", "УП"),
  @("А0274-НСКК. Адм. управление персоналом", "УП"),
  @("А0275- НИСУ. Адм. управление персоналом", "УП"),
  @("А0276-НР. Адм. управление персоналом", "УП"),
  @("А0279-НКиПС. Адм. управление персоналом", "УП"),
  @("А0280-НПИС. Адм. управление персоналом", "УП"),
  @("А0281-НУП. Адм. управление персоналом", "УП"),
  @("А0297-НРБ. Адм. управление персоналом", "УП"),
  @("А0433 - ЦКИС. ГСНиФЭС. Адм.", "УП"),
  @("А0210-НПУ. Адм. управление персоналом", "УП"),
  @("А0288-Адм. время. НИС. Управленческие", "УП"),
  @("А0579 - Адм. управление персоналом Службы сопровождения производственных систем", "УП"),
  @("", ""),
  @("", ""),
  @("А0175-Адм.время НУП. Больничный", "БЛ"),
  @("А0266-Адм. время НИСУ. Больничный", "БЛ"),
  @("А0272-Адм. время НПИС. Больничный", "БЛ"),
  @("А0216-Адм. время НКиПС. Больничный", "БЛ"),
  @("А0268-Адм. время НПУ. Больничный", "БЛ"),
  @("А0262-Адм. время НР. Больничный", "БЛ"),
  @("А0264-Адм. время НСКК. Больничный", "БЛ"),
  @("А0298 -Адм. Время. НРБ. Больничный", "БЛ"),
  @("А0525 - Административное время ООД. Больничный", "БЛ"),
  @("А0603 - Адм. Время. RCBI. Больничный", "БЛ"),
  @("", ""),
  @("", ""),
  @("А0483 - Рабочее время в пути, связанное с командировками и т.п. (происходящие в рабочее время)", "КомВПути"),
  @("", ""),
  @("", ""),
  @("А0536 - ДИСУ. Обучение (повышение квалификации) с отрывом от производства для нужд производства.", "О+Р"),
  @("А0282 - НСКК. Обучение", "О+Р"),
  @("А0316-Адм. Время. НР. Развитие", "О+Р"),
  @("А0536 - ДИСУ. Обучение (повышение квалификации) с отрывом от производства для нужд производства. ", "О+Р"),
  @("А0537 - ДКиПС. Обучение (повышение квалификации) с отрывом от производства для нужд производства.", "О+Р"),
  @("А0629 - Корпоративное обучение не связанное с конкретным ФН, Департаментом или Направлением (для нужд Общества)", "О+Р"),
  @("", ""),
  @("", ""),
  @("В0329- Документооборот (отправка корреспонденции)", "ДП"),
  @("В0519 - НПИС.ФГ.Делопроизводство", "ДП"),
  @("В0261-НР. ФГ Делопроизводство", "ДП"),
  @("В0263-НСКК. ФГ Делопроизводство", "ДП"),
  @("В0332-ДРЭО. Делопроизводство", "ДП"),
  @("В0265-НИСУ. ФГ Делопроизводство", "ДП"),
  @("В0220-НКиПС. ФГ Делопроизводство", "ДП"),
  @("В0267-НПУ. ФГ Делопроизводство", "ДП"),
  @("В0273-НПИС.ФГ Делопроизводство", "ДП"),
  @("", ""),
  @("", ""),
  @("В0305-ФГ. Взаимодействие с кадровыми службами и бухгалтерией", "ФГ"),
  @("В0522 - НПИС.ФГ.Профессионального развития персонала", "ФГ"),
  @("В0520 - НПИС.ФГ.ПЭГНПИС.ФГ.ПЭГ", "ФГ"),
  @("", ""),
  @("", ""),
  @("В0213-Пресейл НИСУ", "Прс"),
  @("В0294-НРБ. Управление процессами проработки коммерческих инициатив и подготовки ТКП", "ФГ"),
  @("В0317-НРБ. Организация взаимодействия с Заказчиками", "ФГ"),
  @("В0292-НУП. Процесс управления сервисами", "ФГ"),
  @("В0530 - ООД. ФГ документооборота (отправка корреспонденции)", "ДП"),
  @("В0217-НКиПС. Пресейл ЦСПБиЭ", "Прс"),
  @("В0277-НКиПС. Пресейл ЦПР", "Прс"),
  @("В0221-НКиПС. Пресейл ЦКИС", "Прс"),
  @("В0278-НКиПС. Пресейл ОСУДиНСИ", "Прс"),
  @("В0381-Поддержка и сопровождение системы ресурсного планирования и учета трудозатрат", "ФГ"),
  @("В0330-Мероприятия по продвижению собственных продуктов и услуг -DES.MES", "Прс"),
  @("В0291-НУП. Информационно-технологическое обслуживание", "ФГ"),
  @("В0219-Пресейл НСКК", "Прс"),
  @("В0377-ФГ.Разработка ПП и ДИ", "ФГ"),
  @("В0315- НРБ. Процесс экономического планирования и контроля ДРЭО", "ФГ"),
  @("В0501 - ЛИС+АМ РИТЭК", "АМ"),
  @("В0502 - ЛИС+АМ ЮГРАНП", "АМ"),
  @("В0511 - ЛИС+АМ МЦПБ", "АМ"),
  @("В0296 -НРБ. ФГ Делопроизводство", "ДП"),
  @("В0512 - ЛИС+АМ ТУРСУНТ", "АМ"),
  @("В0513 - ЛИС+АМ ЭнИнж", "АМ"),
  @("В0514 - ЛИС+АМ Л-ЭПУ", "АМ"),
  @("В0515 - ЛИС+АМ Л-ЧМ", "АМ"),
  @("В0320-НУП. Процесс ресурсного планирования", "ФГ"),
  @("В0533 - ООД. ФГ Служба Заказчика", "ФГ"),
  @("В0534 - ООД. ФГ ведения договоров", "ФГ"),
  @("В0529 - ООД. ФГ Делопроизводство централизованное (кадровое)", "ДП"),
  @("В0531 - ООД. ФГ административно-хозяйственного обеспечения", "ФГ"),
  @("В0532 - ООД. ФГ закупок", "ФГ"),
  @("В0535 - ООД. ФГ ИТО", "ФГ"),
  @("В0528 - ООД. ФГ Делопроизводство централизованное (договорное)", "ДП"),
  @("В0538 - СУП. Процессы формирования и контроля исполнения производственных программ", "ФГ"),
  @("В0539 - СУП. Процессы, обеспечивающие функционирование систем охраны труда и пожарной безопасности.", "ФГ"),
  @("В0527 - ООД. ФГ Делопроизводство централизованное (служебное)", "ДП"),
  @("", ""),
  @("", ""),
  @("В0413 - СРБ. Информационная поддержка коммерческой деятельности", "ФГ"),
  @("В-0365 НОД. Процесс обеспечения деятельности", "ФГ"),
  @("В0293-НУП. Процесс управления портфелем проекта", "ФГ"),
  @("В0300-НУП. Процессы ведения договоров", "ФГ"),
  @("В0542 - ООД. ФГ по проведению инвентаризации", "ФГ"),
  @("В0189- СУП. ФГ Служба Заказчика", "ФГ"),
  @("В0215-НКиПС. Пресейл ЦПСГиД", "Прс"),
  @("В0191-НУП. ФГ Делопроизводство", "ДП"),
  @("В0486 - Каталог работ и услуг Л1", "ФГ"),
  @("В0487 - Количественные показатели работ и услуг (метрики) Л1", "ФГ"),
  @("В0335 - Автоматизация процесса обработки документов в КИС ЭХД", "ФГ"),
  @("В0391 -Пресейл НПУ", "Прс"),
  @("В0507 - ЛИС+АМ НВН", "АМ"),
  @("В0503 - ЛИС+АМ Л-Инж-Пермь", "АМ"),
  @("В0508 - ЛИС+АМ Л-Пермь", "АМ"),
  @("В0510 - ЛИС+АМ УралОйл", "АМ"),
  @("В0333-НОД. ФГ Делопроизводство", "ДП"),
  @("", ""),
  @("", ""),
  @("", ""),
  @("", ""),
  @("", ""),
  @("", ""),
  @("", ""),
  @("", ""),
  @("", ""),
  @("", ""),
  @("", ""),
  @("", ""),
  @("", ""),
  @("", "")
)

for ($i = 0; $i -lt $rows.Count; $i++) {
  $r = $i + 1
  $pair = $rows[$i]
  $ws.Cells.Item($r, 1).Value = $pair[0]
  $ws.Cells.Item($r, 2).Value = $pair[1]
}

$ws.Range("A1:B1").Value2 | Out-Null
Write-Output "rows written"
